$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 9 label (August through date changed from 08-21 to 08-22)
$ws.Range("A9").Value = "August (through 08-22)"

# Update row 9 data values
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 51
$ws.Range("D9").Value = 59
$ws.Range("E9").Value = 39
$ws.Range("F9").Value = 29
$ws.Range("G9").Value = 133
$ws.Range("H9").Value = 109

# Update row 10 (Total) data values
$ws.Range("B10").Value = 187
$ws.Range("C10").Value = 353
$ws.Range("D10").Value = 524
$ws.Range("E10").Value = 464
$ws.Range("F10").Value = 333
$ws.Range("G10").Value = 754
$ws.Range("H10").Value = 1023
